# Apply the recorded navigation/selection changes plus the data edits on
# the "Repayment Schedule" sheet (row 2), then make that sheet the active
# tab, matching the committed state of the workbook.

$wb = $excel.ActiveWorkbook

$wsInput    = $wb.Worksheets.Item("NewLoanInput")
$wsSummary  = $wb.Worksheets.Item("Summary")
$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# --- Data edits on "Repayment Schedule", row 2 -----------------------------
# Add formatting-only (empty) cells in columns B, F and O, copying the
# plain style used across the row.
$wsSchedule.Range("A2").Copy()
$wsSchedule.Range("B2").PasteSpecial(-4122)
$wsSchedule.Range("F2").PasteSpecial(-4122)
$wsSchedule.Range("O2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Clear the values that used to be 0 while keeping their formatting.
$wsSchedule.Range("H2").Value = ""
$wsSchedule.Range("J2").Value = ""

# Set the values that are now populated with 0.
$wsSchedule.Range("I2").Value = 0
$wsSchedule.Range("L2").Value = 0

# --- Selections on the other sheets (no tab switch) -------------------------
$wsInput.Range("B3:B15").Select()
$wsSummary.Range("D4").Select()

# --- Final selection + activate "Repayment Schedule" so it is the active tab
$wsSchedule.Range("D9").Select()
$wsSchedule.Activate()
